$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Hydrogen): corrected Iron & steel demand value, Chemicals
#     (hydrogen) figure removed/corrected to blank ---
$ws.Range("B3").Value = 3649231.461175587
# D3 previously held a number; it becomes a blank text cell (use the
# classic text-prefix trick so it is stored as an empty string rather
# than simply deleting the cell), then reset the style so no stray
# quote-prefix formatting sticks around.
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# --- Row 4 (Methanol): corrected Chemicals figure ---
$ws.Range("C4").Value = 10.00086238711431

# --- Row 5 (Ammonia): corrected Chemicals figure ---
$ws.Range("C5").Value = 1782.208927283131

# --- Row 7: relabel "Other" -> "Biogas" and correct its value ---
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 577.2655423823744

# --- New row 8: "Other" (the row that used to be row 7) with its own
#     corrected value. Copy A7's style (bold/centered/bordered) onto A8. ---
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 540.7476529456285
